$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "Check for first repeated character"
$ws.Range("C17").Value = "Hash - BitSet"
$ws.Range("H17").Value = "CheckFirstDuplicateUsingBitSet"

$ws.Range("A18").Select()
